# Weekly update: insert a new record for "Ají" at Terminal Hortofrutícola
# Agro Chillán as the first data row (row 19, right after the existing
# data block in rows 2-18), shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at row 19; this pushes old rows 19-68 down to 20-69
# and Excel automatically copies the formatting (incl. the date style on
# column D) from the row being pushed down.
$ws.Rows.Item(19).Insert()

# Populate the newly inserted row 19 with the new weekly data point.
$ws.Cells.Item(19, 1).Value  = 7
$ws.Cells.Item(19, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(19, 3).Value  = "Ñuble"
$ws.Cells.Item(19, 4).Value  = 44623
$ws.Cells.Item(19, 5).Value  = 16
$ws.Cells.Item(19, 6).Value  = 100112021
$ws.Cells.Item(19, 7).Value  = "Ají"
$ws.Cells.Item(19, 8).Value  = "Americana (o)"
$ws.Cells.Item(19, 9).Value  = "Primera"
$ws.Cells.Item(19, 10).Value = 60
$ws.Cells.Item(19, 11).Value = 9000
$ws.Cells.Item(19, 12).Value = 9000
$ws.Cells.Item(19, 13).Value = 9000
$ws.Cells.Item(19, 14).Value = "$/caja 15 kilos"
$ws.Cells.Item(19, 15).Value = "Región del Maule"
$ws.Cells.Item(19, 16).Value = 600
$ws.Cells.Item(19, 17).Value = 15
$ws.Cells.Item(19, 18).Value = "Hortaliza"
